# daily auto push: 2026-01-13 09:38 UTC
# A new day's row of data is appended into the sequential log by inserting a
# fresh row right before the first "2026/12/29" entry (row 642) and shifting
# every following row down by one. The new row holds the latest reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 642; rows 642:683 (and the sheet dimension)
# shift down to 643:684 automatically.
$ws.Rows.Item(642).Insert()

# Populate the newly inserted row. The date column is stored as plain text
# in this sheet (not a real date), so prefix with an apostrophe to stop
# Excel's autodetection from coercing it into a date serial number, then
# reapply the Normal style so no stray "quote prefix" number format sticks
# to the cell (keeping it identical in shape to its sibling rows).
$ws.Range("A642").Value = "'2026/01/13"
$ws.Range("A642").Style = "Normal"
$ws.Range("B642").Value = "火"
$ws.Range("C642").Value = 17
$ws.Range("D642").Value = 201
